# Auto-generated Excel COM-interop script
# Applies updated currentAveragePrice market data to Leve profit sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$updates = @{
    "H19" = 2055.125
    "J19" = 3378.6
    "L19" = 3378.6
    "N19" = -3728.6
    "H43" = 6498.5
    "I43" = 2999.6667
    "K43" = 2999.6667
    "M43" = -2930.6667
    "H70" = 3255.2222
    "I70" = 2383.3333
    "K70" = 7149.999899999999
    "M70" = -6879.999899999999
    "H73" = 3255.2222
    "I73" = 2383.3333
    "K73" = 7149.999899999999
    "M73" = -6213.999899999999
    "H88" = 3283.64
    "J88" = 3283.64
    "L88" = 3283.64
    "N88" = -4095.64
    "H91" = 3283.64
    "J91" = 3283.64
    "L91" = 3283.64
    "N91" = -6091.639999999999
}
foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

$ws = $wb.Worksheets.Item("ARM")
$updates = @{
    "H61" = 2504.878
    "I61" = 1655.6364
    "K61" = 1655.6364
    "M61" = -1443.6364
    "H109" = 22295.666
    "J109" = 22295.666
    "L109" = 22295.666
    "N109" = -25069.666
    "H132" = 3428.4285
    "I132" = 0
    "J132" = 3428.4285
    "K132" = 0
    "L132" = 10285.2855
    "N132" = -15345.2855
    "H136" = 2504.878
    "I136" = 1655.6364
    "K136" = 4966.9092
    "M136" = -2416.9092
}
foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
$ws.Range("M132").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$updates = @{
    "H86" = 2158.8
    "I86" = 1949.75
    "K86" = 1949.75
    "M86" = -826.75
    "H89" = 2158.8
    "I89" = 1949.75
    "K89" = 9748.75
    "M89" = -4132.75
    "H122" = 0
    "J122" = 0
    "L122" = 0
}
foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$updates = @{
    "H31" = 1502.9678
    "I31" = 1117.5
    "J31" = 2445.2222
    "K31" = 1117.5
    "L31" = 2445.2222
    "M31" = -822.5
    "N31" = -3035.2222
    "H34" = 1502.9678
    "I34" = 1117.5
    "J34" = 2445.2222
    "K34" = 1117.5
    "L34" = 2445.2222
    "M34" = -915.5
    "N34" = -2849.2222
    "H58" = 4502.048
    "I58" = 1963.4546
    "K58" = 1963.4546
    "M58" = -1760.4546
    "H99" = 2090.2163
    "I99" = 2033.5
    "J99" = 2733
    "K99" = 2033.5
    "L99" = 2733
    "M99" = -535.5
    "N99" = -5729
    "H107" = 861.1667
    "I107" = 315.85715
    "J107" = 992.7931
    "K107" = 315.85715
    "L107" = 992.7931
    "M107" = 1604.14285
    "N107" = -4832.7931
    "H126" = 2090.2163
    "I126" = 2033.5
    "J126" = 2733
    "K126" = 6100.5
    "L126" = 8199
    "M126" = -3630.5
    "N126" = -13139
    "H136" = 4502.048
    "I136" = 1963.4546
    "K136" = 5890.3638
    "M136" = -3340.3638
}
foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

$ws = $wb.Worksheets.Item("GSM")
$updates = @{
    "H102" = 2727.2354
    "I102" = 2729.1667
    "K102" = 2729.1667
    "M102" = -1107.1667
    "H126" = 5421782.5
    "I126" = 4266.8945
    "K126" = 12800.6835
    "M126" = -10330.6835
}
foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

$ws = $wb.Worksheets.Item("LTW")
$updates = @{
    "H7" = 4673.433
    "I7" = 4665.2666
    "K7" = 4665.2666
    "M7" = -4553.2666
    "H40" = 4060.2104
    "I40" = 3557.647
    "K40" = 3557.647
    "M40" = -3421.647
    "H46" = 2033
    "J46" = 3000
    "L46" = 3000
    "N46" = -3376
    "H82" = 2243.4443
    "J82" = 2135.875
    "L82" = 2135.875
    "N82" = -2857.875
    "H85" = 2243.4443
    "J85" = 2135.875
    "L85" = 2135.875
    "N85" = -4631.875
    "H110" = 74633
    "J110" = 74633
    "L110" = 74633
    "N110" = -82813
    "H114" = 39832.332
    "J114" = 39832.332
    "L114" = 39832.332
    "N114" = -48510.332
    "H122" = 5600.826
    "I122" = 5518.875
    "J122" = 5788.143
    "K122" = 16556.625
    "L122" = 17364.429
    "M122" = -14106.625
    "N122" = -22264.429
    "H124" = 50214.5
    "J124" = 50214.5
    "L124" = 50214.5
    "N124" = -60034.5
    "H126" = 4673.433
    "I126" = 4665.2666
    "K126" = 13995.7998
    "M126" = -11525.7998
    "H132" = 4261.5654
    "I132" = 4064.9412
    "J132" = 4818.6665
    "K132" = 12194.8236
    "L132" = 14455.9995
    "M132" = -9664.8236
    "N132" = -19515.9995
    "H136" = 1775.6923
    "I136" = 1446.52
    "K136" = 4339.559999999999
    "M136" = -1789.559999999999
}
foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

$ws = $wb.Worksheets.Item("WVR")
$updates = @{
    "H111" = 50644
    "J111" = 50644
    "L111" = 50644
    "N111" = -58824
    "H131" = 79995
    "J131" = 79995
    "L131" = 79995
    "N131" = -90075
    "H132" = 2149.9092
    "I132" = 1990
    "K132" = 5970
    "M132" = -3440
}
foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
